# edit.ps1 - applies the "Updated cryptos list" diff to cryptos.xlsx
# Columns: A=index(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so that numeric-looking values
# (e.g. "1.001", "0.3920", "1.000") are stored verbatim instead of being
# auto-parsed/normalized into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (D) and Volume 1h (E) updates for rows 2-34 ---
$ws.Range("D2").Value = "29.269.07"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.901.09"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "326.05"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.3920"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "0.07894"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").Value = "0.9894"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").Value = "22.06"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").Value = "1.908.95"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").Value = "7.086"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "5.755"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "0.06980"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "88.41"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "0.00001001"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "17.09"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "29.261.12"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "5.322"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "11.08"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "2.093"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "156.47"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "6.010"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").Value = "118.66"
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "1.922"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("D30").Value = "0.09379"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "0.9077"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "5.296"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").Value = "1.329"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "3.224"
$ws.Range("E34").Value = "  -1.68%  "

# --- Rows 35 and 36 swapped places (Hedera now ranks above TrustWalletToken) ---
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.05801"
$ws.Range("E35").Value = "  -0.98%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.181"
$ws.Range("E36").Value = "  +1.66%  "

# --- Price (D) and Volume 1h (E) updates for rows 37-51 ---
$ws.Range("D37").Value = "0.02091"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "1.000"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "7.773"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").Value = "0.5715"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "0.1787"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("D42").Value = "9.761"
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("D43").Value = "12.04"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "2.212"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").Value = "0.5359"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").Value = "0.07053"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "1.862"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "2.586"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "113.21"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "1.066"
$ws.Range("E50").Value = "  -4.65%  "
$ws.Range("D51").Value = "71.36"
$ws.Range("E51").Value = "  -0.70%  "

# Reset the style on column D back to "Normal" (General) so that we do not
# leave a dangling explicit Text-number-format style on these cells; the
# values themselves remain text because they have already been committed.
$ws.Range("D2:D51").Style = "Normal"
